# Updates from 0409 meeting - add hospital capacity variables to the
# "Relevant Variable" sheet, remove the stray "in R" note in E1, and
# widen column D to fit the new (longer) note text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Relevant Variable")

# Remove the old E1 "in R" annotation (no longer relevant / moved out).
$ws.Range("E1").ClearContents()

# New rows describing the hospital-capacity variables pulled from the
# state-wise bed/ventilator estimate PDFs (2020-2031 in the source note
# naming, one row per metric).
$hospitalVars = @(
    "hostpital_public",
    "hospital_private",
    "hospital_total",
    "hospital_beds_public",
    "hospital_beds_private",
    "hospital_beds_total",
    "icu_beds_public",
    "icu_beds_private",
    "icu_beds_total",
    "ventilators_public",
    "ventilators_private",
    "ventilators_total"
)

$startYear = 2020
$startRow = 13

# Fill column-by-column (A, then B, then C, then D) so new shared-string
# entries are interned in the same order the author produced them in
# (type all variable names, then fill dataset/level, then the notes).
for ($i = 0; $i -lt $hospitalVars.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $hospitalVars[$i]
}

for ($i = 0; $i -lt $hospitalVars.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = "hostpital_capacity"
}

for ($i = 0; $i -lt $hospitalVars.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = "state"
}

for ($i = 0; $i -lt $hospitalVars.Length; $i++) {
    $row = $startRow + $i
    $year = $startYear + $i
    $note = "from 'State-wise-estimates-of-current-beds-and-ventilators_24Apr$year.pdf'"
    $ws.Cells.Item($row, 4).Value = $note
}

# Widen column D to fit the longer note text added above (matches the
# ~70.6-character-wide column the author ended up with after autofitting
# to the new, longer note strings).
$ws.Columns.Item(4).ColumnWidth = 69.8

# Leave the active selection on E1, matching the author's final cursor
# position after clearing that cell.
$ws.Range("E1").Select()

Write-Output "edit applied"
